# Daily refresh of the "剩余" (days remaining) tracker.
# For every data row: remaining (col E) counts down by 1 day.
# When a row's remaining count hits 1 (i.e. about to expire), the
# item is considered replenished: remaining is reset back to the
# total day count (col D) and the start date (col F) is advanced by
# that same number of days.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $total = $ws.Cells.Item($r, 4).Value2      # D: 总天 (total days)
    $remaining = $ws.Cells.Item($r, 5).Value2  # E: 剩余 (days remaining)
    $startRaw = $ws.Cells.Item($r, 6).Value2   # F: 开始时间 (start date, yyyymmdd)

    if ($null -eq $total -or $null -eq $remaining -or $null -eq $startRaw) {
        continue
    }

    $startText = [string]$startRaw
    if ($startText.Length -ne 8) {
        # Not a well-formed yyyymmdd value (data issue) - leave the row untouched.
        continue
    }

    $year = [int]$startText.Substring(0, 4)
    $month = [int]$startText.Substring(4, 2)
    $day = [int]$startText.Substring(6, 2)

    $startDate = $null
    try {
        $startDate = Get-Date -Year $year -Month $month -Day $day
    } catch {
        $startDate = $null
    }

    if ($null -eq $startDate) {
        continue
    }

    if ($remaining -le 1) {
        # Item restocked: reset the countdown and roll the start date forward.
        $newRemaining = $total
        $newStartDate = $startDate.AddDays($total)
        $newStart = [int]$newStartDate.ToString("yyyyMMdd")

        $ws.Cells.Item($r, 5).Value = $newRemaining
        $ws.Cells.Item($r, 6).Value = $newStart
    } else {
        $ws.Cells.Item($r, 5).Value = $remaining - 1
    }
}
